# Scheduled market-data refresh: update price/profit columns (H:N) across
# the Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 132.5
$ws.Range("I6").Value = 132.5
$ws.Range("K6").Value = 397.5
$ws.Range("M6").Value = -285.5

$ws.Range("H18").Value = 389.58334
$ws.Range("I18").Value = 349.5
$ws.Range("J18").Value = 590
$ws.Range("K18").Value = 349.5
$ws.Range("L18").Value = 590
$ws.Range("M18").Value = -65.5
$ws.Range("N18").Value = -1158

$ws.Range("H43").Value = 1242.4445
$ws.Range("I43").Value = 965.6667
$ws.Range("J43").Value = 1380.8334
$ws.Range("K43").Value = 965.6667
$ws.Range("L43").Value = 1380.8334
$ws.Range("M43").Value = -896.6667
$ws.Range("N43").Value = -1518.8334

# Row 69: L69 changes and a new M69 cell is introduced
$ws.Range("H69").Value = 3291.8572
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 3673.8333
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 11021.4999
$ws.Range("M69").Value = -2126
$ws.Range("N69").Value = -12769.4999

# Row 72: same pattern as row 69
$ws.Range("H72").Value = 3291.8572
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 3673.8333
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 33064.4997
$ws.Range("M72").Value = -4632
$ws.Range("N72").Value = -41800.4997

$ws.Range("H116").Value = 2010.6552
$ws.Range("I116").Value = 1829.8823
$ws.Range("J116").Value = 2266.75
$ws.Range("K116").Value = 1829.8823
$ws.Range("L116").Value = 2266.75
$ws.Range("M116").Value = 1612.1177
$ws.Range("N116").Value = -9150.75

$ws.Range("H132").Value = 4016.5
$ws.Range("I132").Value = 2284.2666
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 6852.7998
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -4322.7998
$ws.Range("N132").Value = -95060

$ws.Range("H135").Value = 78949130
$ws.Range("I135").Value = 31251696
$ws.Range("J135").Value = 333335420
$ws.Range("K135").Value = 281265264
$ws.Range("L135").Value = 3000018780
$ws.Range("M135").Value = -281262729
$ws.Range("N135").Value = -3000023850

$ws.Range("H138").Value = 8998.661
$ws.Range("J138").Value = 15768.613
$ws.Range("L138").Value = 47305.839
$ws.Range("N138").Value = -57585.839

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 19999.143
$ws.Range("J56").Value = 19999.143
$ws.Range("L56").Value = 19999.143
$ws.Range("N56").Value = -21483.143

# Row 106: values collapse to 0 and the N106 cell is removed entirely
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 50991.875
$ws.Range("J140").Value = 50991.875
$ws.Range("L140").Value = 50991.875
$ws.Range("N140").Value = -61351.875

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 500000000
$ws.Range("I6").Value = 500000000
$ws.Range("K6").Value = 500000000
$ws.Range("M6").Value = -499999887

$ws.Range("H7").Value = 206.83333
$ws.Range("I7").Value = 48.2
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 48.2
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 64.8
$ws.Range("N7").Value = -1226

$ws.Range("H31").Value = 12523
$ws.Range("I31").Value = 10923.8
$ws.Range("J31").Value = 15188.333
$ws.Range("K31").Value = 10923.8
$ws.Range("L31").Value = 15188.333
$ws.Range("M31").Value = -10628.8
$ws.Range("N31").Value = -15778.333

$ws.Range("H34").Value = 12523
$ws.Range("I34").Value = 10923.8
$ws.Range("J34").Value = 15188.333
$ws.Range("K34").Value = 10923.8
$ws.Range("L34").Value = 15188.333
$ws.Range("M34").Value = -10721.8
$ws.Range("N34").Value = -15592.333

$ws.Range("H50").Value = 26218.643
$ws.Range("J50").Value = 26218.643
$ws.Range("L50").Value = 26218.643
$ws.Range("N50").Value = -27468.643

$ws.Range("H51").Value = 25630.572
$ws.Range("J51").Value = 25630.572
$ws.Range("L51").Value = 25630.572
$ws.Range("N51").Value = -27102.572

$ws.Range("H59").Value = 32675.666
$ws.Range("J59").Value = 32675.666
$ws.Range("L59").Value = 32675.666
$ws.Range("N59").Value = -34965.666

$ws.Range("H60").Value = 8776.368
$ws.Range("J60").Value = 9152.833000000001
$ws.Range("L60").Value = 9152.833000000001
$ws.Range("N60").Value = -10174.833

$ws.Range("H61").Value = 25630.572
$ws.Range("J61").Value = 25630.572
$ws.Range("L61").Value = 25630.572
$ws.Range("N61").Value = -26326.572

$ws.Range("H62").Value = 3834.25
$ws.Range("I62").Value = 3399.9
$ws.Range("J62").Value = 6006
$ws.Range("K62").Value = 3399.9
$ws.Range("L62").Value = 6006
$ws.Range("M62").Value = -2775.9
$ws.Range("N62").Value = -7254

$ws.Range("H65").Value = 3834.25
$ws.Range("I65").Value = 3399.9
$ws.Range("J65").Value = 6006
$ws.Range("K65").Value = 16999.5
$ws.Range("L65").Value = 30030
$ws.Range("M65").Value = -13879.5
$ws.Range("N65").Value = -36270

# Row 74: values collapse to 0 and the N74 cell is removed entirely
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""

# Row 77: values collapse to 0 and the N77 cell is removed entirely
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 260.16666
$ws.Range("I6").Value = 260.16666
$ws.Range("K6").Value = 780.4999799999999
$ws.Range("M6").Value = -667.4999799999999

$ws.Range("H14").Value = 1525.05
$ws.Range("I14").Value = 1525.05
$ws.Range("K14").Value = 4575.15
$ws.Range("M14").Value = -4402.15

# Row 62: values collapse to 0 and the N62 cell is removed entirely
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

# Row 65: values collapse to 0 and the N65 cell is removed entirely
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13750
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 17000
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 17000
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -18996

$ws.Range("H83").Value = 13750
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 17000
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 85000
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -94984

$ws.Range("H107").Value = 566.7778
$ws.Range("I107").Value = 150.5
$ws.Range("J107").Value = 899.8
$ws.Range("K107").Value = 150.5
$ws.Range("L107").Value = 899.8
$ws.Range("M107").Value = 1769.5
$ws.Range("N107").Value = -4739.8

$ws.Range("H132").Value = 7101.615
$ws.Range("I132").Value = 9926.538
$ws.Range("K132").Value = 29779.614
$ws.Range("M132").Value = -27249.614

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 9375
$ws.Range("I56").Value = 5833.3335
$ws.Range("J56").Value = 20000
$ws.Range("K56").Value = 5833.3335
$ws.Range("L56").Value = 20000
$ws.Range("M56").Value = -5142.3335
$ws.Range("N56").Value = -21382

$ws.Range("H68").Value = 3899
$ws.Range("I68").Value = 3899
$ws.Range("K68").Value = 3899
$ws.Range("M68").Value = -3150

$ws.Range("H71").Value = 3899
$ws.Range("I71").Value = 3899
$ws.Range("K71").Value = 19495
$ws.Range("M71").Value = -15751

$ws.Range("H82").Value = 2259.4
$ws.Range("I82").Value = 2124.25
$ws.Range("K82").Value = 2124.25
$ws.Range("M82").Value = -1763.25

$ws.Range("H85").Value = 2259.4
$ws.Range("I85").Value = 2124.25
$ws.Range("K85").Value = 2124.25
$ws.Range("M85").Value = -876.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 7: values change and a new N7 cell is introduced
$ws.Range("H7").Value = 1804.5
$ws.Range("J7").Value = 1804.5
$ws.Range("L7").Value = 1804.5
$ws.Range("N7").Value = -2030.5

# Row 14: values change and a new N14 cell is introduced
$ws.Range("H14").Value = 851666.7
$ws.Range("I14").Value = 1275000
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 1275000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -1274832
$ws.Range("N14").Value = -5336

# Row 42: M42 cell is removed and a new N42 cell is introduced
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 10000
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -10756

$ws.Range("H58").Value = 13616.6
$ws.Range("I58").Value = 10027.667
$ws.Range("K58").Value = 10027.667
$ws.Range("M58").Value = -9719.666999999999

$ws.Range("H81").Value = 28574670
$ws.Range("J81").Value = 33336784
$ws.Range("L81").Value = 66673568
$ws.Range("N81").Value = -66675690

$ws.Range("H84").Value = 28574670
$ws.Range("J84").Value = 33336784
$ws.Range("L84").Value = 333367840
$ws.Range("N84").Value = -333378448

$ws.Range("H136").Value = 3481.4656
$ws.Range("I136").Value = 2980.6924
$ws.Range("J136").Value = 3888.3438
$ws.Range("K136").Value = 8942.0772
$ws.Range("L136").Value = 11665.0314
$ws.Range("M136").Value = -6392.0772
$ws.Range("N136").Value = -16765.0314
